# Apply strikethrough formatting to the "Break point visualization" paragraph
# (paragraph mark + run), matching the other "done" items in the list that
# already use <w:strike/> on both the paragraph rPr and the run rPr.

$d = $word.ActiveDocument

$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "Break point visualization") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    # Applying StrikeThrough to the whole paragraph Range (including the
    # paragraph mark) sets <w:strike/> in both <w:pPr><w:rPr> and the run's
    # <w:rPr>, exactly like the sibling "strike" items elsewhere in the doc.
    $target.Range.Font.StrikeThrough = 1
}
